$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57. This shifts the existing rows 57-111
# down to 58-112 (their values/formatting move with them), matching the
# diff where old row 57 data reappears (unchanged) as new row 58, ...,
# and old row 111 data reappears (unchanged) as new row 112.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new data record.
$ws.Range("A57").Value = 7
$ws.Range("B57").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C57").Value = "Ñuble"
$ws.Range("D57").Value = 45159
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = 100112001
$ws.Range("G57").Value = "Berenjena"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 80
$ws.Range("K57").Value = 8000
$ws.Range("L57").Value = 8000
$ws.Range("M57").Value = 8000
$ws.Range("N57").Value = "`$/caja 60 unidades"
$ws.Range("O57").Value = "Región de Arica y Parinacota"
$ws.Range("P57").Value = 133
$ws.Range("Q57").Value = 60
$ws.Range("R57").Value = "Hortaliza"
